$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated experimental values (column S, "hclZORASR") for rows 3-7 ---
$ws.Range("S3").Value = 193.36
$ws.Range("S4").Value = 193.41
$ws.Range("S5").Value = 193.42
$ws.Range("S6").Value = 193.2
$ws.Range("S7").Value = 193.22

# --- X7 ("hclX2CAMF2s") now carries its own explicit formula instead of
#     continuing the shared formula used by the rows above it ---
$ws.Range("X7").Formula = "=W7*27.211"

# --- View state: scrolled so column F is the leftmost visible column,
#     with the new selection sitting on Y8 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Y8").Select()

# --- Workbook window geometry (size/position) ---
$win = $wb.Windows.Item(1)
$win.Left = 1620
$win.Top = 6580
$win.Width = 11700
$win.Height = 11560
